$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (column D) values - force text storage so values like
# "149.80" or "1.00" are not reinterpreted as numbers (matches original
# inline-string formatting with no numeric coercion / trailing zeros kept).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.488.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.483.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "491.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.492.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0982"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.912.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.433.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.504.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.411"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.162"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.597.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0796"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.873"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0557"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.611"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "266.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0928"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.895.42"
$ws.Range("D51").Style = "Normal"

# Update coin name / link (columns B & C) and volume-change % (column E)
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("E6").Value = "  +7.97%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("E10").Value = "  +5.94%  "
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("E16").Value = "  +3.55%  "
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  +4.25%  "
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +3.96%  "
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("E38").Value = "  +4.09%  "
$ws.Range("E39").Value = "  +5.74%  "
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("E41").Value = "  +3.47%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  +6.03%  "
$ws.Range("E46").Value = "  +9.02%  "
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("E51").Value = "  -3.58%  "
